$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-7 (Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster, and the 16 NATMI metric columns E:T), following the
# "Natmi following Dr Hou advice" re-run that expands the single ligand/
# receptor pair onto every Sending x Target cluster combination
# (ECs, FAPs, sCs) x (FAPs, sCs).
$rows = @(
    ,@("ECs", "Gdf9", "Bmpr1b", "FAPs", 2, 0.6666666666666666, 0.5498363333333334, 1.649509, 0.1730423789436426, 0.1730423789436426, 3, 1, 1.292736333333333, 3.878209, 0.9024289107934966, 0.9024289107934966, 0.7107934054867778, 6.397140649381001, 0.1561584455512269, 0.1561584455512269)
    ,@("ECs", "Gdf9", "Bmpr1b", "sCs", 2, 0.6666666666666666, 0.5498363333333334, 1.649509, 0.1730423789436426, 0.1730423789436426, 2, 0.6666666666666666, 0.1397713333333333, 0.419314, 0.09757108920650338, 0.09757108920650338, 0.07685135742511111, 0.691662216826, 0.01688393339241572, 0.01688393339241572)
    ,@("FAPs", "Gdf9", "Bmpr1b", "FAPs", 3, 1, 2.071011333333333, 6.213034, 0.6517807322165176, 0.6517807322165177, 3, 1, 1.292736333333333, 3.878209, 0.9024289107934966, 0.9024289107934966, 2.677271597345111, 24.095444376106, 0.5881857762503396, 0.5881857762503397)
    ,@("FAPs", "Gdf9", "Bmpr1b", "sCs", 3, 1, 2.071011333333333, 6.213034, 0.6517807322165176, 0.6517807322165177, 2, 0.6666666666666666, 0.1397713333333333, 0.419314, 0.09757108920650338, 0.09757108920650338, 0.2894680154084444, 2.605212138676, 0.06359495596617792, 0.06359495596617794)
    ,@("sCs", "Gdf9", "Bmpr1b", "FAPs", 2, 0.6666666666666666, 0.5566186666666667, 1.669856, 0.1751768888398398, 0.1751768888398398, 3, 1, 1.292736333333333, 3.878209, 0.9024289107934966, 0.9024289107934966, 0.7195611742115555, 6.476050567904, 0.1580846889919301, 0.15808468899193)
    ,@("sCs", "Gdf9", "Bmpr1b", "sCs", 2, 0.6666666666666666, 0.5566186666666667, 1.669856, 0.1751768888398398, 0.1751768888398398, 2, 0.6666666666666666, 0.1397713333333333, 0.419314, 0.09757108920650338, 0.09757108920650338, 0.07779933319822223, 0.700193998784, 0.01709219984790974, 0.01709219984790973)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
